$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037747584065066
$ws.Range("D2").Value = 1.044889327199931
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.053922114303413
$ws.Range("I2").Value = 1.040941350211184
$ws.Range("J2").Value = 1.042848772806245
$ws.Range("K2").Value = 1.04765910864655
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.056666747259589
$ws.Range("N2").Value = 1.044329737680056
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038529238543601
$ws.Range("D3").Value = 1.045487605081989
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.054656576556018
$ws.Range("I3").Value = 1.041109093685107
$ws.Range("J3").Value = 1.043275777547245
$ws.Range("K3").Value = 1.048069271492459
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.057214551684232
$ws.Range("N3").Value = 1.044757348816767
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03903564103574
$ws.Range("D4").Value = 1.045875297316095
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.055132741521048
$ws.Range("I4").Value = 1.041216797564238
$ws.Range("J4").Value = 1.043552004578557
$ws.Range("K4").Value = 1.048334542272475
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.05756926528707
$ws.Range("N4").Value = 1.045033968122147
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039248679141596
$ws.Range("D5").Value = 1.046038416831656
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.055333139099935
$ws.Range("I5").Value = 1.041261875044098
$ws.Range("J5").Value = 1.043668111464576
$ws.Range("K5").Value = 1.048446029130776
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.057718444429575
$ws.Range("N5").Value = 1.045150239893261
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039284457713145
$ws.Range("D6").Value = 1.046065813100842
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.055366799446783
$ws.Range("I6").Value = 1.041269431928851
$ws.Range("J6").Value = 1.043687605170943
$ws.Range("K6").Value = 1.048464746291042
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.057743495588663
$ws.Range("N6").Value = 1.045169761282927
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03903848708805
$ws.Range("D7").Value = 1.045877476404094
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.055135418390069
$ws.Range("I7").Value = 1.041217400683382
$ws.Range("J7").Value = 1.043553556080624
$ws.Range("K7").Value = 1.048336032096637
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.057571258402114
$ws.Range("N7").Value = 1.045035521827524
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03801161839754
$ws.Range("D8").Value = 1.045091399581232
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.054170137742433
$ws.Range("I8").Value = 1.04099821282586
$ws.Range("J8").Value = 1.042993095488969
$ws.Range("K8").Value = 1.047797751607642
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.056851828053293
$ws.Range("N8").Value = 1.044474265317544
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036206974431009
$ws.Range("D9").Value = 1.043710654266841
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.052476320770582
$ws.Range("I9").Value = 1.040605596443955
$ws.Range("J9").Value = 1.042004987999624
$ws.Range("K9").Value = 1.046848282372132
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.055586071116457
$ws.Range("N9").Value = 1.043484754602226
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035007237981406
$ws.Range("D10").Value = 1.042793244746187
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.051352031164037
$ws.Range("I10").Value = 1.040339608817338
$ws.Range("J10").Value = 1.04134599069868
$ws.Range("K10").Value = 1.046214746156294
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.054743660777173
$ws.Range("N10").Value = 1.042824821449527
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034488559398557
$ws.Range("D11").Value = 1.042396752730264
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.050866395711892
$ws.Range("I11").Value = 1.040223436444542
$ws.Range("J11").Value = 1.04106059311387
$ws.Range("K11").Value = 1.045940303327842
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.0543792480764
$ws.Range("N11").Value = 1.042539018567411
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034296023472734
$ws.Range("D12").Value = 1.042249592942792
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.050686189971553
$ws.Range("I12").Value = 1.040180135689371
$ws.Range("J12").Value = 1.040954577990594
$ws.Range("K12").Value = 1.045838346800308
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.054243944297401
$ws.Range("N12").Value = 1.0424328528905
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034337317422519
$ws.Range("D13").Value = 1.042281153981163
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.050724836493647
$ws.Range("I13").Value = 1.040189430585905
$ws.Range("J13").Value = 1.040977318832773
$ws.Range("K13").Value = 1.045860217533955
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.05427296488558
$ws.Range("N13").Value = 1.042455626027283
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034472641770176
$ws.Range("D14").Value = 1.042384586098975
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.050851496146375
$ws.Range("I14").Value = 1.040219860229034
$ws.Range("J14").Value = 1.041051829980152
$ws.Range("K14").Value = 1.045931875892663
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.054368062685076
$ws.Range("N14").Value = 1.042530242989038
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034556036107811
$ws.Range("D15").Value = 1.042448329337983
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.050929559374101
$ws.Range("I15").Value = 1.04023858919397
$ws.Range("J15").Value = 1.041097738035846
$ws.Range("K15").Value = 1.045976024862383
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.054426662961431
$ws.Range("N15").Value = 1.042576216239437
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035041678309995
$ws.Range("D16").Value = 1.042819574634863
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.051384286436094
$ws.Range("I16").Value = 1.040347297834721
$ws.Range("J16").Value = 1.041364930714565
$ws.Range("K16").Value = 1.046232957668167
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.054767853309863
$ws.Range("N16").Value = 1.042843788362407
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035346528442338
$ws.Range("D17").Value = 1.043052649819732
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.051669844689988
$ws.Range("I17").Value = 1.040415221250685
$ws.Range("J17").Value = 1.041532521897881
$ws.Range("K17").Value = 1.046394094373347
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.054981969697604
$ws.Range("N17").Value = 1.043011617544426
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035524420894855
$ws.Range("D18").Value = 1.043188671091892
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.05183652059744
$ws.Range("I18").Value = 1.040454743488122
$ws.Range("J18").Value = 1.041630270334042
$ws.Range("K18").Value = 1.046488071259206
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.055106894446678
$ws.Range("N18").Value = 1.043109504794578
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035585090885422
$ws.Range("D19").Value = 1.043235063061182
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.051893372147309
$ws.Range("I19").Value = 1.040468203185038
$ws.Range("J19").Value = 1.04166359919169
$ws.Range("K19").Value = 1.046520112967501
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.055149496319188
$ws.Range("N19").Value = 1.043142880983027
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035313812771872
$ws.Range("D20").Value = 1.043027635549115
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.051639195121404
$ws.Range("I20").Value = 1.040407943670513
$ws.Range("J20").Value = 1.041514541422695
$ws.Range("K20").Value = 1.046376807101779
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.054958993476882
$ws.Range("N20").Value = 1.042993611534903
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034432788673376
$ws.Range("D21").Value = 1.042354124727885
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.050814193024462
$ws.Range("I21").Value = 1.040210903569011
$ws.Range("J21").Value = 1.041029888452209
$ws.Range("K21").Value = 1.045910774729117
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.054340057212064
$ws.Range("N21").Value = 1.042508270301608
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033879574510286
$ws.Range("D22").Value = 1.041931328074651
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.050296529119856
$ws.Range("I22").Value = 1.040086154006259
$ws.Range("J22").Value = 1.040725135286214
$ws.Range("K22").Value = 1.045617668406812
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.053951228370767
$ws.Range("N22").Value = 1.042203084351161
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034172774855004
$ws.Range("D23").Value = 1.042155396763642
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.050570852510318
$ws.Range("I23").Value = 1.040152367619374
$ws.Range("J23").Value = 1.040886693339487
$ws.Range("K23").Value = 1.045773057985862
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.054157322887108
$ws.Range("N23").Value = 1.042364871835401
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035328595332736
$ws.Range("D24").Value = 1.043038938197656
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.051653043988326
$ws.Range("I24").Value = 1.040411232393515
$ws.Range("J24").Value = 1.041522666040068
$ws.Range("K24").Value = 1.046384618511176
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.054969375336333
$ws.Range("N24").Value = 1.043001747690164
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036672934234355
$ws.Range("D25").Value = 1.044067074188096
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.052913354886069
$ws.Range("I25").Value = 1.040707848468223
$ws.Range("J25").Value = 1.042260488979071
$ws.Range("K25").Value = 1.047093847066926
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.055913055478289
$ws.Range("N25").Value = 1.04374061842237
